# Update automática del mapa (2025-08-28 07:06:15)
# Adds a new data row (row 84) to the PEBCOM sheet, mirroring the
# existing rows: text-typed columns for A-H, J-L, O-P and
# number-typed columns for I, M, N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 84

# Columns that must be stored as *text*, even though some values look
# numeric (Caso, OT, Comuna, etc. are plain strings in this sheet).
$textValues = @{
    1  = "7064"
    2  = "8/28/2025"
    3  = "GAONA AV. 1189"
    4  = "6"
    5  = "809257408"
    6  = "PEBCOM"
    7  = "Pendiente"
    8  = "Piden aplomo de columna ver con Pedro si hay que colocar una R400 por tension de la red"
    10 = "Cambio"
    11 = "Sin equipos"
    12 = "Terminal"
    15 = "Almagro"
    16 = "Capital Sur"
}

foreach ($col in $textValues.Keys) {
    $cell = $ws.Cells.Item($newRow, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $textValues[$col]
    $cell.ClearFormats()
}

# Numeric columns: Attachments, Coordenada_X, Coordenada_Y
$ws.Cells.Item($newRow, 9).Value = 1
$ws.Cells.Item($newRow, 13).Value = -58.446008
$ws.Cells.Item($newRow, 14).Value = -34.607602
